# Generate Report for Handoff
# - Refreshes the existing handoff row (old guid -> new guid, new timestamps/hashes)
# - Appends a new row for a newly discovered file that is also ready for handoff

$wb = $excel.ActiveWorkbook

$oldGuidName   = "6796f835-e5bd-410e-bbe2-c2e2278f81f7.md"
$newGuidName   = "740ec174-6d4a-4531-ae7d-5e19eca46094.md"
$newGuid2Name  = "ffffe684983d-7fd2-423b-a493-8c85ed6d4b24.md"

$newZhXlf = "740ec174-6d4a-4531-ae7d-5e19eca46094.76a4fffd6dc20842fbee121cd9345c0f6cdeb40b.zh-cn.xlf"
$newDeXlf = "740ec174-6d4a-4531-ae7d-5e19eca46094.76a4fffd6dc20842fbee121cd9345c0f6cdeb40b.de-de.xlf"

$overviewDate = "2016-11-09 00:56:02"
$zhHandoffDate = "2016-11-09 00:55:48"
$deHandoffDate = "2016-11-09 00:56:02"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a82a8571995b5ff241af1a223dc84f190c65bfeb/e2e/"

# ---------------------------------------------------------------------------
# Sheet "Overview"  (columns: A File Name, B Path And Name, C Extension,
#   D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# -- refresh row 2 for the renamed file --
$wsOverview.Cells.Item(2, 1).Value = $newGuidName
$wsOverview.Cells.Item(2, 2).Hyperlinks.Delete()
$wsOverview.Cells.Item(2, 2).Value = "e2e\" + $newGuidName
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(2, 2), $baseUrl + $newGuidName, "", "", "e2e\" + $newGuidName)
$wsOverview.Cells.Item(2, 2).Style = "HyperLink"
$wsOverview.Cells.Item(2, 7).Value = $overviewDate

# -- add row 3 for the new file --
$loOverview.ListRows.Add() | Out-Null
$wsOverview.Cells.Item(3, 1).Value = $newGuid2Name
$wsOverview.Cells.Item(3, 2).Value = "e2e\" + $newGuid2Name
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), $baseUrl + $newGuid2Name, "", "", "e2e\" + $newGuid2Name)
$wsOverview.Cells.Item(3, 2).Style = "HyperLink"
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = ""
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = $overviewDate

# ---------------------------------------------------------------------------
# Sheet "zh-cn"  (columns: A Source File Name, B File Extension, C Status,
#   D Source Path, E Priority, F Content Duplicate, G Latest Handoff File,
#   H Latest Handoff Datetime, I Latest Target File, J Latest Handback File,
#   K Latest Handback DateTime, L Reference Tokens, M To be localized,
#   N Dependency From, O Has metadata, P Error Detail)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# -- refresh row 2 for the renamed file --
$wsZh.Cells.Item(2, 1).Hyperlinks.Delete()
$wsZh.Cells.Item(2, 1).Value = $newGuidName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 1), $baseUrl + $newGuidName, "", "", $newGuidName)
$wsZh.Cells.Item(2, 1).Style = "HyperLink"
$wsZh.Cells.Item(2, 7).Value = $newZhXlf
$wsZh.Cells.Item(2, 8).Value = $zhHandoffDate

# -- add row 3 for the new file --
$loZh.ListRows.Add() | Out-Null
$wsZh.Cells.Item(3, 1).Value = $newGuid2Name
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 1), $baseUrl + $newGuid2Name, "", "", $newGuid2Name)
$wsZh.Cells.Item(3, 1).Style = "HyperLink"
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "True"
$wsZh.Cells.Item(3, 7).Value = $newZhXlf
$wsZh.Cells.Item(3, 8).Value = $zhHandoffDate
$wsZh.Cells.Item(3, 9).Value = ""
$wsZh.Cells.Item(3, 10).Value = ""
$wsZh.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(3, 12).Value = ""
$wsZh.Cells.Item(3, 13).Value = "True"
$wsZh.Cells.Item(3, 14).Value = ""
$wsZh.Cells.Item(3, 15).Value = "False"
$wsZh.Cells.Item(3, 16).Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de"  (same columns as zh-cn)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# -- refresh row 2 for the renamed file --
$wsDe.Cells.Item(2, 1).Hyperlinks.Delete()
$wsDe.Cells.Item(2, 1).Value = $newGuidName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 1), $baseUrl + $newGuidName, "", "", $newGuidName)
$wsDe.Cells.Item(2, 1).Style = "HyperLink"
$wsDe.Cells.Item(2, 7).Value = $newDeXlf
$wsDe.Cells.Item(2, 8).Value = $deHandoffDate

# -- add row 3 for the new file --
$loDe.ListRows.Add() | Out-Null
$wsDe.Cells.Item(3, 1).Value = $newGuid2Name
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 1), $baseUrl + $newGuid2Name, "", "", $newGuid2Name)
$wsDe.Cells.Item(3, 1).Style = "HyperLink"
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "True"
$wsDe.Cells.Item(3, 7).Value = $newDeXlf
$wsDe.Cells.Item(3, 8).Value = $deHandoffDate
$wsDe.Cells.Item(3, 9).Value = ""
$wsDe.Cells.Item(3, 10).Value = ""
$wsDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(3, 12).Value = ""
$wsDe.Cells.Item(3, 13).Value = "True"
$wsDe.Cells.Item(3, 14).Value = ""
$wsDe.Cells.Item(3, 15).Value = "False"
$wsDe.Cells.Item(3, 16).Value = ""

Write-Host "Report regenerated for handoff."
